# Generate Report for Handoff
# Update status text and timestamps across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (B2, C2) and handoff date (D2)
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-27-14 06:27:48"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (E2)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-14 06:27:46"

# de-de sheet: Status (C2) and Latest Handoff Datetime (E2)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-14 06:27:48"
